# Added passport priority import handler to hca.py
# Reflects the corresponding new reference row added to the
# PassportPriority lookup sheet: id=5, passport_priority="Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PassportPriority")

# Make this the active sheet (mirrors the author switching to it in Excel).
$ws.Activate()

# Append the new lookup row right after the existing data (rows 1-5).
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Done"

# Leave the selection on the newly added cell, as in the saved workbook.
$ws.Range("B6").Select()
